$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chain Lightning (row 13): Buff shock damage from 8 to 8(11), set Code column to "x".
$ws.Range("C13").Value = "Deal 8(11) Shock damage to ALL enemies."
$ws.Range("F13").Value = "x"

# Fireball (row 8): Reduce energy cost from 3 to 2, and weaken self-damage from 6(7) to 6.
$ws.Range("C8").Value = "Deal 6 Fire damage to yourself. Deal 36(42) Fire damage."
$ws.Range("D8").Value = 2

# Update the active selection to reflect where the edit was made.
$ws.Range("C8").Select()
